# Normalize the "Recorded By" column (G): the list of recorders in each
# cell is a comma-separated string; re-sort the names within each cell
# using an ordinal (case-sensitive, ASCII) alphabetical order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $list = New-Object System.Collections.Generic.List[string]
            foreach ($p in $parts) { [void]$list.Add($p) }
            $list.Sort([System.StringComparer]::Ordinal)
            $newVal = [string]::Join(", ", $list)
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
